$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 6.407607666666666
$ws.Range("H2").Value = 19.222823
$ws.Range("I2").Value = 0.01049006948643881
$ws.Range("J2").Value = 0.01049006948643881
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 0.6186053333333333
$ws.Range("N2").Value = 1.855816
$ws.Range("O2").Value = 0.0556943868446899
$ws.Range("P2").Value = 0.0556943868446899
$ws.Range("Q2").Value = 3.963780276507555
$ws.Range("R2").Value = 35.67402248856799
$ws.Range("S2").Value = 0.0005842379880054003
$ws.Range("T2").Value = 0.0005842379880054003
$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 6.407607666666666
$ws.Range("H3").Value = 19.222823
$ws.Range("I3").Value = 0.01049006948643881
$ws.Range("J3").Value = 0.01049006948643881
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 2.338622
$ws.Range("N3").Value = 7.015866
$ws.Range("O3").Value = 0.2105512373287584
$ws.Range("P3").Value = 0.2105512373287584
$ws.Range("Q3").Value = 14.98497225663533
$ws.Range("R3").Value = 134.864750309718
$ws.Range("S3").Value = 0.002208697110034344
$ws.Range("T3").Value = 0.002208697110034344
$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 6.407607666666666
$ws.Range("H4").Value = 19.222823
$ws.Range("I4").Value = 0.01049006948643881
$ws.Range("J4").Value = 0.01049006948643881
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 0.165314
$ws.Range("N4").Value = 0.495942
$ws.Range("O4").Value = 0.01488357983794147
$ws.Range("P4").Value = 0.01488357983794148
$ws.Range("Q4").Value = 1.059267253807333
$ws.Range("R4").Value = 9.533405284265999
$ws.Range("S4").Value = 0.0001561297867069657
$ws.Range("T4").Value = 0.0001561297867069657
$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 6.407607666666666
$ws.Range("H5").Value = 19.222823
$ws.Range("I5").Value = 0.01049006948643881
$ws.Range("J5").Value = 0.01049006948643881
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 7.984598333333333
$ws.Range("N5").Value = 23.953795
$ws.Range("O5").Value = 0.7188707959886103
$ws.Range("P5").Value = 0.7188707959886103
$ws.Range("Q5").Value = 51.16217349592056
$ws.Range("R5").Value = 460.459561463285
$ws.Range("S5").Value = 0.007541004601692096
$ws.Range("T5").Value = 0.007541004601692096
$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 572.1502276666666
$ws.Range("H6").Value = 1716.450683
$ws.Range("I6").Value = 0.9366827616690507
$ws.Range("J6").Value = 0.9366827616690508
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 0.6186053333333333
$ws.Range("N6").Value = 1.855816
$ws.Range("O6").Value = 0.0556943868446899
$ws.Range("P6").Value = 0.0556943868446899
$ws.Range("Q6").Value = 353.9351823024809
$ws.Range("R6").Value = 3185.416640722328
$ws.Range("S6").Value = 0.05216797207914858
$ws.Range("T6").Value = 0.05216797207914859
$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 572.1502276666666
$ws.Range("H7").Value = 1716.450683
$ws.Range("I7").Value = 0.9366827616690507
$ws.Range("J7").Value = 0.9366827616690508
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 2.338622
$ws.Range("N7").Value = 7.015866
$ws.Range("O7").Value = 0.2105512373287584
$ws.Range("P7").Value = 0.2105512373287584
$ws.Range("Q7").Value = 1338.043109726275
$ws.Range("R7").Value = 12042.38798753648
$ws.Range("S7").Value = 0.1972197144539372
$ws.Range("T7").Value = 0.1972197144539372
$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 572.1502276666666
$ws.Range("H8").Value = 1716.450683
$ws.Range("I8").Value = 0.9366827616690507
$ws.Range("J8").Value = 0.9366827616690508
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 0.165314
$ws.Range("N8").Value = 0.495942
$ws.Range("O8").Value = 0.01488357983794147
$ws.Range("P8").Value = 0.01488357983794148
$ws.Range("Q8").Value = 94.58444273648732
$ws.Range("R8").Value = 851.259984628386
$ws.Range("S8").Value = 0.01394119266612482
$ws.Range("T8").Value = 0.01394119266612483
$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 572.1502276666666
$ws.Range("H9").Value = 1716.450683
$ws.Range("I9").Value = 0.9366827616690507
$ws.Range("J9").Value = 0.9366827616690508
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 7.984598333333333
$ws.Range("N9").Value = 23.953795
$ws.Range("O9").Value = 0.7188707959886103
$ws.Range("P9").Value = 0.7188707959886103
$ws.Range("Q9").Value = 4568.389754243553
$ws.Range("R9").Value = 41115.50778819199
$ws.Range("S9").Value = 0.6733538824698402
$ws.Range("T9").Value = 0.6733538824698403
$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 0.356025
$ws.Range("H10").Value = 1.068075
$ws.Range("I10").Value = 0.0005828582496300428
$ws.Range("J10").Value = 0.0005828582496300428
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 0.6186053333333333
$ws.Range("N10").Value = 1.855816
$ws.Range("O10").Value = 0.0556943868446899
$ws.Range("P10").Value = 0.0556943868446899
$ws.Range("Q10").Value = 0.2202389638
$ws.Range("R10").Value = 1.9821506742
$ws.Range("S10").Value = 0.00003246193283051444
$ws.Range("T10").Value = 0.00003246193283051444
$ws.Range("E11").Value = 3
$ws.Range("G11").Value = 0.356025
$ws.Range("H11").Value = 1.068075
$ws.Range("I11").Value = 0.0005828582496300428
$ws.Range("J11").Value = 0.0005828582496300428
$ws.Range("K11").Value = 3
$ws.Range("M11").Value = 2.338622
$ws.Range("N11").Value = 7.015866
$ws.Range("O11").Value = 0.2105512373287584
$ws.Range("P11").Value = 0.2105512373287584
$ws.Range("Q11").Value = 0.83260789755
$ws.Range("R11").Value = 7.493471077949999
$ws.Range("S11").Value = 0.0001227215256468798
$ws.Range("T11").Value = 0.0001227215256468798
$ws.Range("E12").Value = 3
$ws.Range("G12").Value = 0.356025
$ws.Range("H12").Value = 1.068075
$ws.Range("I12").Value = 0.0005828582496300428
$ws.Range("J12").Value = 0.0005828582496300428
$ws.Range("K12").Value = 3
$ws.Range("M12").Value = 0.165314
$ws.Range("N12").Value = 0.495942
$ws.Range("O12").Value = 0.01488357983794147
$ws.Range("P12").Value = 0.01488357983794148
$ws.Range("Q12").Value = 0.05885591685
$ws.Range("R12").Value = 0.52970325165
$ws.Range("S12").Value = 0.000008675017292571564
$ws.Range("T12").Value = 0.000008675017292571564
$ws.Range("E13").Value = 3
$ws.Range("G13").Value = 0.356025
$ws.Range("H13").Value = 1.068075
$ws.Range("I13").Value = 0.0005828582496300428
$ws.Range("J13").Value = 0.0005828582496300428
$ws.Range("K13").Value = 3
$ws.Range("M13").Value = 7.984598333333333
$ws.Range("N13").Value = 23.953795
$ws.Range("O13").Value = 0.7188707959886103
$ws.Range("P13").Value = 0.7188707959886103
$ws.Range("Q13").Value = 2.842716621625
$ws.Range("R13").Value = 25.584449594625
$ws.Range("S13").Value = 0.000418999773860077
$ws.Range("T13").Value = 0.000418999773860077
$ws.Range("E14").Value = 3
$ws.Range("G14").Value = 31.91218566666667
$ws.Range("H14").Value = 95.736557
$ws.Range("I14").Value = 0.05224431059488034
$ws.Range("J14").Value = 0.05224431059488035
$ws.Range("K14").Value = 3
$ws.Range("M14").Value = 0.6186053333333333
$ws.Range("N14").Value = 1.855816
$ws.Range("O14").Value = 0.0556943868446899
$ws.Range("P14").Value = 0.0556943868446899
$ws.Range("Q14").Value = 19.74104825172356
$ws.Range("R14").Value = 177.669434265512
$ws.Range("S14").Value = 0.002909714844705397
$ws.Range("T14").Value = 0.002909714844705397
$ws.Range("E15").Value = 3
$ws.Range("G15").Value = 31.91218566666667
$ws.Range("H15").Value = 95.736557
$ws.Range("I15").Value = 0.05224431059488034
$ws.Range("J15").Value = 0.05224431059488035
$ws.Range("K15").Value = 3
$ws.Range("M15").Value = 2.338622
$ws.Range("N15").Value = 7.015866
$ws.Range("O15").Value = 0.2105512373287584
$ws.Range("P15").Value = 0.2105512373287584
$ws.Range("Q15").Value = 74.63053946815134
$ws.Range("R15").Value = 671.674855213362
$ws.Range("S15").Value = 0.01100010423914002
$ws.Range("T15").Value = 0.01100010423914002
$ws.Range("E16").Value = 3
$ws.Range("G16").Value = 31.91218566666667
$ws.Range("H16").Value = 95.736557
$ws.Range("I16").Value = 0.05224431059488034
$ws.Range("J16").Value = 0.05224431059488035
$ws.Range("K16").Value = 3
$ws.Range("M16").Value = 0.165314
$ws.Range("N16").Value = 0.495942
$ws.Range("O16").Value = 0.01488357983794147
$ws.Range("P16").Value = 0.01488357983794148
$ws.Range("Q16").Value = 5.275531061299334
$ws.Range("R16").Value = 47.479779551694
$ws.Range("S16").Value = 0.0007775823678171132
$ws.Range("T16").Value = 0.0007775823678171134
$ws.Range("E17").Value = 3
$ws.Range("G17").Value = 31.91218566666667
$ws.Range("H17").Value = 95.736557
$ws.Range("I17").Value = 0.05224431059488034
$ws.Range("J17").Value = 0.05224431059488035
$ws.Range("K17").Value = 3
$ws.Range("M17").Value = 7.984598333333333
$ws.Range("N17").Value = 23.953795
$ws.Range("O17").Value = 0.7188707959886103
$ws.Range("P17").Value = 0.7188707959886103
$ws.Range("Q17").Value = 254.8059844870906
$ws.Range("R17").Value = 2293.253860383815
$ws.Range("S17").Value = 0.03755690914321781
$ws.Range("T17").Value = 0.03755690914321782
